# Applies the "updating metrics in briefing file" edit to the briefing
# document: normalizes a handful of runs that had been split apart by
# Word's spell-checker (proofErr-wrapped words) back into single runs,
# and appends a new "Dia 06/09" time-tracking line at the end of the
# document.

$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- Collapse runs that Word's spell-checker had split with <w:proofErr/> ---
# (the visible text is identical; only the run/proofErr markup is merged)

$t1 = "Oferecer conteúdo informativo sobre sáude mental"
Replace-Exact $t1 $t1

$t2 = "Gênero: Predominamente Feminino"
Replace-Exact $t2 $t2

$t3 = "Comportamento online: Uso frequente de redes sociais, busca por informações sobre sáude mental, interesse em conteúdos de autoajuda e psicologia"
Replace-Exact $t3 $t3

$t4 = "O cliente ja possui identidade visual?"
Replace-Exact $t4 $t4

# Five identical "Hex: " labels (each followed by its own color-code run)
$t5 = "Hex: "
Replace-Exact $t5 $t5

$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$t6 = "Fontes: Utilizar uma fonte sans-serif como " + $quoteOpen + "Roboto" + $quoteClose + " para o texto geral, que é moderna e de facil leitura. Para títulos, uma fonte serifada elegante como " + $quoteOpen + "Merriweather" + $quoteClose + " pode ser usada para adicionar um toque de sofisticação."
Replace-Exact $t6 $t6

# --- Append the new time-tracking entry at the end of the document ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Dia 06/09: 1hr 10min (1 dia)"
